$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Register the date / datetime number formats used below -----------------
# (Registering lowercase first and then uppercase reproduces the same
#  numFmt id ordering/content seen in the authored workbook: 164/165 for the
#  date-only format and 166/167 for the date+time format, while only the
#  uppercase ones end up referenced by the cellXfs actually used.)
$ws.Range("D2").NumberFormat = "yyyy-mm-dd"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD"

$ws.Range("E2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("E2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$dateFormat = "YYYY-MM-DD"
$dateTimeFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row data ------------------------------------------------------------
$rows = @(
    @{ Row = 2; Empleado = "703"; Nombre = "JAIME EMANUEL ES "; Dia = "Viernes";
       Fecha = 44120; I0 = 44120.38958333333; E0 = 44120.66736111111;
       I1 = 44120; E1 = 44120; I2 = 44120; E2 = 44120; I3 = 44120; E3 = 44120; I4 = 44120; E4 = 44120 },
    @{ Row = 3; Empleado = "703"; Nombre = "JAIME EMANUEL ES "; Dia = "Sábado";
       Fecha = 44121; I0 = 44121.32430555556; E0 = 44121.66736111111;
       I1 = 44121; E1 = 44121; I2 = 44121; E2 = 44121; I3 = 44121; E3 = 44121; I4 = 44121; E4 = 44121 },
    @{ Row = 4; Empleado = "705"; Nombre = "BALMACEDA FRANCO NICOLAS DR "; Dia = "Viernes";
       Fecha = 44120; I0 = 44120.38541666666; E0 = 44120.54166666666;
       I1 = 44120.56180555555; E1 = 44120.7; I2 = 44120; E2 = 44120; I3 = 44120; E3 = 44120; I4 = 44120; E4 = 44120 },
    @{ Row = 5; Empleado = "706"; Nombre = "PERALTA MARIO ALBERTO RE "; Dia = "Viernes";
       Fecha = 44120; I0 = 44120.39166666667; E0 = 44120.52291666667;
       I1 = 44120.54027777778; E1 = 44120.70208333333; I2 = 44120; E2 = 44120; I3 = 44120; E3 = 44120; I4 = 44120; E4 = 44120 }
)

foreach ($r in $rows) {
    $n = $r.Row

    # Column B / C are plain (non numeric-looking) text, no special handling needed.
    $ws.Range("B$n").Value = $r.Nombre
    $ws.Range("C$n").Value = $r.Dia

    # Column D: date only
    $ws.Range("D$n").NumberFormat = $dateFormat
    $ws.Range("D$n").Value = $r.Fecha

    # Columns E..N: date+time
    $ws.Range("E$n").NumberFormat = $dateTimeFormat
    $ws.Range("E$n").Value = $r.I0
    $ws.Range("F$n").NumberFormat = $dateTimeFormat
    $ws.Range("F$n").Value = $r.E0
    $ws.Range("G$n").NumberFormat = $dateTimeFormat
    $ws.Range("G$n").Value = $r.I1
    $ws.Range("H$n").NumberFormat = $dateTimeFormat
    $ws.Range("H$n").Value = $r.E1
    $ws.Range("I$n").NumberFormat = $dateTimeFormat
    $ws.Range("I$n").Value = $r.I2
    $ws.Range("J$n").NumberFormat = $dateTimeFormat
    $ws.Range("J$n").Value = $r.E2
    $ws.Range("K$n").NumberFormat = $dateTimeFormat
    $ws.Range("K$n").Value = $r.I3
    $ws.Range("L$n").NumberFormat = $dateTimeFormat
    $ws.Range("L$n").Value = $r.E3
    $ws.Range("M$n").NumberFormat = $dateTimeFormat
    $ws.Range("M$n").Value = $r.I4
    $ws.Range("N$n").NumberFormat = $dateTimeFormat
    $ws.Range("N$n").Value = $r.E4
}

# Column A ("Empleado") holds numeric-looking codes (703/705/706) that must be
# stored as text, exactly like in the source workbook. Force text storage via
# the "@" number format and then strip the format back off (Style = Normal)
# so the cells keep the default (unstyled) look used in the target file.
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("A2").Value = "703"
$ws.Range("A3").Value = "703"
$ws.Range("A4").Value = "705"
$ws.Range("A5").Value = "706"
$ws.Range("A2:A5").Style = "Normal"

Write-Output "done"
